$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18: clear the "MEC-1NA" entry from column E (move it out)
$ws.Range("E18").Value = "-"

# Row 19: clear column E, add the "MEC-1NA" entry to column F
$ws.Range("E19").Value = "-"
$ws.Range("F19").Value = "['MEC-1NA-Tec. Mat. Não Metal.', -, -, -]"

# Row 20: add the "MEC-1NA" entry to column B, clear column E, add it to column F
$ws.Range("B20").Value = "['MEC-1NA-Tec. Mat. Não Metal.', -, -, -]"
$ws.Range("E20").Value = "-"
$ws.Range("F20").Value = "['MEC-1NA-Tec. Mat. Não Metal.', -, -, -]"

# Row 21: clear column E, add the "MEC-1NA" entry to column F
$ws.Range("E21").Value = "-"
$ws.Range("F21").Value = "['MEC-1NA-Tec. Mat. Não Metal.', -, -, -]"
